$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for 1a5419d9-... (row 3 and row 4 share the same value)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-03 08:19:19"
$wsOverview.Range("G4").Value = "2016-09-03 08:19:19"

# Sheet "zh-cn": Status ht -> mt, Correspond Handoff/Handback Datetime updates for 1a5419d9-... rows (3 and 4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-03 08:19:14"
$wsZhCn.Range("H4").Value = "2016-09-03 08:19:14"
$wsZhCn.Range("K3").Value = "2016-09-03 08:19:31"
$wsZhCn.Range("K4").Value = "2016-09-03 08:19:31"

# Sheet "de-de": Status ht -> mt, Correspond Handoff Datetime (shared with Overview's G column) and Handback Datetime updates
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-03 08:19:19"
$wsDeDe.Range("H4").Value = "2016-09-03 08:19:19"
$wsDeDe.Range("K3").Value = "2016-09-03 08:19:38"
$wsDeDe.Range("K4").Value = "2016-09-03 08:19:38"
